$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NSW press conference held 25/05/2023 -> update "Last Date" for the NSW row
$ws.Range("B5").Value2 = 45071

# Update the NSW "News Link" cell text to the new YouTube link
$ws.Range("C5").Value = "https://www.youtube.com/watch?v=AOISAo2T3Rw"

# Final cursor position left on B6 by the author
$ws.Range("B6").Select()
